# Scheduled-runner refresh of market-board derived profit figures.
# Updates currentAveragePrice* / Leve/Craft price & profit columns (H..N)
# on a handful of rows across several job sheets, per the latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 19624.273
$ws.Range("I34").Value = 19624.273
$ws.Range("K34").Value = 19624.273
$ws.Range("M34").Value = -19421.273

# Row 36
$ws.Range("H36").Value = 19624.273
$ws.Range("I36").Value = 19624.273
$ws.Range("K36").Value = 19624.273
$ws.Range("M36").Value = -18909.273

# Row 51
$ws.Range("H51").Value = 42220.2
$ws.Range("I51").Value = 101250.5
$ws.Range("J51").Value = 2866.6667
$ws.Range("K51").Value = 101250.5
$ws.Range("L51").Value = 2866.6667
$ws.Range("M51").Value = -100766.5
$ws.Range("N51").Value = -3834.6667

# Row 134
$ws.Range("H134").Value = 76780
$ws.Range("J134").Value = 76780
$ws.Range("L134").Value = 76780
$ws.Range("N134").Value = -86920

# Row 135
$ws.Range("H135").Value = 714.17645
$ws.Range("I135").Value = 714.17645
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6427.58805
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3892.58805
$ws.Range("N135").ClearContents()

# Row 137
$ws.Range("H137").Value = 1872.6
$ws.Range("I137").Value = 1872.6
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5617.799999999999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3067.799999999999
$ws.Range("N137").ClearContents()

# Row 138
$ws.Range("H138").Value = 4187.875
$ws.Range("I138").Value = 2136.5
$ws.Range("J138").Value = 5120.3184
$ws.Range("K138").Value = 6409.5
$ws.Range("L138").Value = 15360.9552
$ws.Range("M138").Value = -1269.5
$ws.Range("N138").Value = -25640.9552

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13013.173
$ws.Range("I32").Value = 13948.192
$ws.Range("K32").Value = 13948.192
$ws.Range("M32").Value = -13661.192

# Row 123
$ws.Range("H123").Value = 30426.125
$ws.Range("J123").Value = 30426.125
$ws.Range("L123").Value = 30426.125
$ws.Range("N123").Value = -40226.125

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 4420.9165
$ws.Range("I22").Value = 4777.364
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 4777.364
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -4604.364
$ws.Range("N22").Value = -846

# Row 86
$ws.Range("H86").Value = 128500.5
$ws.Range("I86").Value = 4400.4
$ws.Range("J86").Value = 335334
$ws.Range("K86").Value = 4400.4
$ws.Range("L86").Value = 335334
$ws.Range("M86").Value = -3277.4
$ws.Range("N86").Value = -337580

# Row 89
$ws.Range("H89").Value = 128500.5
$ws.Range("I89").Value = 4400.4
$ws.Range("J89").Value = 335334
$ws.Range("K89").Value = 22002
$ws.Range("L89").Value = 1676670
$ws.Range("M89").Value = -16386
$ws.Range("N89").Value = -1687902

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2007.2368
$ws.Range("I31").Value = 1164.8438
$ws.Range("K31").Value = 1164.8438
$ws.Range("M31").Value = -869.8438000000001

# Row 34
$ws.Range("H34").Value = 2007.2368
$ws.Range("I34").Value = 1164.8438
$ws.Range("K34").Value = 1164.8438
$ws.Range("M34").Value = -962.8438000000001

# Row 68
$ws.Range("H68").Value = 32000
$ws.Range("J68").Value = 32000
$ws.Range("L68").Value = 32000
$ws.Range("N68").Value = -33498

# Row 71
$ws.Range("H71").Value = 32000
$ws.Range("J71").Value = 32000
$ws.Range("L71").Value = 96000
$ws.Range("N71").Value = -103488

# Row 87
$ws.Range("H87").Value = 41876.668
$ws.Range("J87").Value = 41876.668
$ws.Range("L87").Value = 41876.668
$ws.Range("N87").Value = -44248.668

# Row 90
$ws.Range("H90").Value = 41876.668
$ws.Range("J90").Value = 41876.668
$ws.Range("L90").Value = 125630.004
$ws.Range("N90").Value = -137486.004

# Row 99
$ws.Range("H99").Value = 2131.4
$ws.Range("I99").Value = 2155.5557
$ws.Range("J99").Value = 1914
$ws.Range("K99").Value = 2155.5557
$ws.Range("L99").Value = 1914
$ws.Range("M99").Value = -657.5556999999999
$ws.Range("N99").Value = -4910

# Row 122
$ws.Range("H122").Value = 2919
$ws.Range("I122").Value = 2887.3333
$ws.Range("J122").Value = 3014
$ws.Range("K122").Value = 8661.999899999999
$ws.Range("L122").Value = 9042
$ws.Range("M122").Value = -6211.999899999999
$ws.Range("N122").Value = -13942

# Row 126
$ws.Range("H126").Value = 2131.4
$ws.Range("I126").Value = 2155.5557
$ws.Range("J126").Value = 1914
$ws.Range("K126").Value = 6466.6671
$ws.Range("L126").Value = 5742
$ws.Range("M126").Value = -3996.6671
$ws.Range("N126").Value = -10682

$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 679.8
$ws.Range("J44").Value = 649.75
$ws.Range("L44").Value = 1949.25
$ws.Range("N44").Value = -2745.25

# Row 117
$ws.Range("H117").Value = 50736
$ws.Range("J117").Value = 63170
$ws.Range("L117").Value = 189510
$ws.Range("N117").Value = -196394

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 263.75
$ws.Range("I107").Value = 326
$ws.Range("J107").Value = 201.5
$ws.Range("K107").Value = 326
$ws.Range("L107").Value = 201.5
$ws.Range("M107").Value = 1594
$ws.Range("N107").Value = -4041.5

# Row 122
$ws.Range("H122").Value = 2685.9333
$ws.Range("I122").Value = 1781.3
$ws.Range("K122").Value = 5343.9
$ws.Range("M122").Value = -2893.9

# Row 123
$ws.Range("H123").Value = 9325.25
$ws.Range("J123").Value = 9325.25
$ws.Range("L123").Value = 9325.25
$ws.Range("N123").Value = -14225.25

# Row 131
$ws.Range("H131").Value = 48205.332
$ws.Range("J131").Value = 48205.332
$ws.Range("L131").Value = 48205.332
$ws.Range("N131").Value = -58285.332

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4431.8887
$ws.Range("I40").Value = 4816.4
$ws.Range("J40").Value = 3951.25
$ws.Range("K40").Value = 4816.4
$ws.Range("L40").Value = 3951.25
$ws.Range("M40").Value = -4680.4
$ws.Range("N40").Value = -4223.25

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 2964
$ws.Range("I62").Value = 2964
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2964
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2340
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 2964
$ws.Range("I65").Value = 2964
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 14820
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11700
$ws.Range("N65").ClearContents()

# Row 113
$ws.Range("H113").Value = 617.45
$ws.Range("I113").Value = 505.3
$ws.Range("J113").Value = 729.6
$ws.Range("K113").Value = 1515.9
$ws.Range("L113").Value = 2188.8
$ws.Range("M113").Value = 654.0999999999999
$ws.Range("N113").Value = -6528.8

# Row 122
$ws.Range("H122").Value = 3422.4443
$ws.Range("I122").Value = 2182.75
$ws.Range("J122").Value = 4414.2
$ws.Range("K122").Value = 6548.25
$ws.Range("L122").Value = 13242.6
$ws.Range("M122").Value = -4098.25
$ws.Range("N122").Value = -18142.6

# Row 123
$ws.Range("H123").Value = 40087.69
$ws.Range("J123").Value = 40087.69
$ws.Range("L123").Value = 40087.69
$ws.Range("N123").Value = -49887.69

# Row 125
$ws.Range("H125").Value = 60715
$ws.Range("J125").Value = 60715
$ws.Range("L125").Value = 60715
$ws.Range("N125").Value = -70555

# Row 132
$ws.Range("H132").Value = 3705.4211
$ws.Range("I132").Value = 3125
$ws.Range("K132").Value = 9375
$ws.Range("M132").Value = -6845
